$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.104.30"
$ws.Range("E2").Value = "'  +0.34%  "
$ws.Range("D3").Value = "'2.672.12"
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("D5").Value = "'578.70"
$ws.Range("E5").Value = "'  +0.26%  "
$ws.Range("D6").Value = "'145.03"
$ws.Range("E6").Value = "'  +1.36%  "
$ws.Range("E7").Value = "'  +0.15%  "
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("D9").Value = "'6.62"
$ws.Range("E9").Value = "'  +1.39%  "
$ws.Range("E10").Value = "'  +1.10%  "
$ws.Range("D11").Value = "'0.386"
$ws.Range("E11").Value = "'  +5.45%  "
$ws.Range("E12").Value = "'  +0.93%  "
$ws.Range("D13").Value = "'3.144.23"
$ws.Range("E13").Value = "'  +2.27%  "
$ws.Range("D14").Value = "'25.98"
$ws.Range("E14").Value = "'  +11.53%  "
$ws.Range("D15").Value = "'61.108.39"
$ws.Range("E15").Value = "'  +0.45%  "
$ws.Range("E16").Value = "'  +1.37%  "
$ws.Range("D17").Value = "'2.671.65"
$ws.Range("E17").Value = "'  +1.86%  "
$ws.Range("E18").Value = "'  +3.11%  "
$ws.Range("E19").Value = "'  +2.64%  "
$ws.Range("D20").Value = "'351.99"
$ws.Range("E20").Value = "'  +0.88%  "
$ws.Range("E21").Value = "'  +0.16%  "
$ws.Range("E22").Value = "'  +0.30%  "
$ws.Range("D23").Value = "'0.530"
$ws.Range("E23").Value = "'  +1.39%  "
$ws.Range("D24").Value = "'64.02"
$ws.Range("E24").Value = "'  +1.02%  "
$ws.Range("E25").Value = "'  +0.27%  "
$ws.Range("D26").Value = "'0.162"
$ws.Range("E26").Value = "'  +0.75%  "
$ws.Range("D27").Value = "'8.15"
$ws.Range("E27").Value = "'  +5.00%  "
$ws.Range("D28").Value = "'1.97"
$ws.Range("E28").Value = "'  +5.99%  "
$ws.Range("D29").Value = "'0.0₃0818"
$ws.Range("E29").Value = "'  +3.20%  "
$ws.Range("D30").Value = "'6.92"
$ws.Range("E30").Value = "'  +8.50%  "
$ws.Range("D32").Value = "'165.39"
$ws.Range("E32").Value = "'  +2.07%  "
$ws.Range("E33").Value = "'  +2.14%  "
$ws.Range("E34").Value = "'  +11.36%  "
$ws.Range("D35").Value = "'4.50"
$ws.Range("E35").Value = "'  +6.38%  "
$ws.Range("E36").Value = "'  +6.33%  "
$ws.Range("E37").Value = "'  +4.58%  "
$ws.Range("D38").Value = "'336.72"
$ws.Range("E38").Value = "'  +11.79%  "
$ws.Range("D39").Value = "'4.05"
$ws.Range("E39").Value = "'  +4.47%  "
$ws.Range("D40").Value = "'38.63"
$ws.Range("E40").Value = "'  +1.94%  "
$ws.Range("D41").Value = "'0.889"
$ws.Range("E41").Value = "'  +5.21%  "
$ws.Range("D42").Value = "'5.20"
$ws.Range("E42").Value = "'  +4.72%  "
$ws.Range("D43").Value = "'20.43"
$ws.Range("E43").Value = "'  +2.56%  "
$ws.Range("D44").Value = "'134.39"
$ws.Range("E44").Value = "'  -0.40%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.100"
$ws.Range("E45").Value = "'  +1.55%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0563"
$ws.Range("E46").Value = "'  +2.78%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'20.66"
$ws.Range("E47").Value = "'  +3.84%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0249"
$ws.Range("E48").Value = "'  +3.10%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.616"
$ws.Range("E49").Value = "'  +1.46%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "'  +0.11%  "
$ws.Range("D51").Value = "'2.097.90"
$ws.Range("E51").Value = "'  +3.42%  "
